{"js": "// The document contains one 20x5 table of simple arithmetic expressions\n// (\"NN+NN=\" / \"NN-NN=\"). The edit replaces the text of every cell with a\n// new expression, keeping the table shape (20 rows x 5 columns) and all\n// per-run formatting (font/size) untouched. Because a couple of \"before\"\n// values repeat (e.g. \"58+16=\" appears twice but maps to two different\n// \"after\" values), the replacement must be positional (row-major order),\n// not a text search/replace.\n\nconst newValues = [\n  [\"88-41=\", \"8+28=\", \"43+36=\", \"53-8=\", \"4-4=\"],\n  [\"28+6=\", \"47-13=\", \"63+3=\", \"88-41=\", \"65-40=\"],\n  [\"44+3=\", \"97-77=\", \"72-40=\", \"72+8=\", \"19+62=\"],\n  [\"21+70=\", \"10+25=\", \"70-10=\", \"92-62=\", \"99-24=\"],\n  [\"15+37=\", \"32+12=\", \"60-54=\", \"28+35=\", \"37-33=\"],\n  [\"14+79=\", \"50+30=\", \"42-26=\", \"15+42=\", \"29+53=\"],\n  [\"96-4=\", \"20+40=\", \"14-12=\", \"30-14=\", \"59-53=\"],\n  [\"38+17=\", \"20+17=\", \"20+64=\", \"2+19=\", \"99-24=\"],\n  [\"86-60=\", \"2+20=\", \"91-0=\", \"94-75=\", \"72-36=\"],\n  [\"6+20=\", \"16+69=\", \"78-66=\", \"72-18=\", \"31-11=\"],\n  [\"14+11=\", \"20+33=\", \"13+42=\", \"91-72=\", \"98-68=\"],\n  [\"73+11=\", \"83+1=\", \"39+50=\", \"34+53=\", \"22+13=\"],\n  [\"40+51=\", \"48-12=\", \"13+40=\", \"30+11=\", \"22+7=\"],\n  [\"65+17=\", \"58-3=\", \"59+1=\", \"20+69=\", \"64-18=\"],\n  [\"52+13=\", \"25+0=\", \"94-25=\", \"64-38=\", \"57-19=\"],\n  [\"86-6=\", \"42-9=\", \"30+63=\", \"73-22=\", \"15-7=\"],\n  [\"55-45=\", \"54-27=\", \"68+26=\", \"36+13=\", \"13-12=\"],\n  [\"62+6=\", \"72+19=\", \"88-86=\", \"13+54=\", \"97-4=\"],\n  [\"57-0=\", \"22+52=\", \"83+11=\", \"60-31=\", \"14-5=\"],\n  [\"20+37=\", \"24+31=\", \"67-65=\", \"43-3=\", \"59-6=\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.values = newValues;\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$newValues = @(\n    @(\"88-41=\", \"8+28=\", \"43+36=\", \"53-8=\", \"4-4=\"),\n    @(\"28+6=\", \"47-13=\", \"63+3=\", \"88-41=\", \"65-40=\"),\n    @(\"44+3=\", \"97-77=\", \"72-40=\", \"72+8=\", \"19+62=\"),\n    @(\"21+70=\", \"10+25=\", \"70-10=\", \"92-62=\", \"99-24=\"),\n    @(\"15+37=\", \"32+12=\", \"60-54=\", \"28+35=\", \"37-33=\"),\n    @(\"14+79=\", \"50+30=\", \"42-26=\", \"15+42=\", \"29+53=\"),\n    @(\"96-4=\", \"20+40=\", \"14-12=\", \"30-14=\", \"59-53=\"),\n    @(\"38+17=\", \"20+17=\", \"20+64=\", \"2+19=\", \"99-24=\"),\n    @(\"86-60=\", \"2+20=\", \"91-0=\", \"94-75=\", \"72-36=\"),\n    @(\"6+20=\", \"16+69=\", \"78-66=\", \"72-18=\", \"31-11=\"),\n    @(\"14+11=\", \"20+33=\", \"13+42=\", \"91-72=\", \"98-68=\"),\n    @(\"73+11=\", \"83+1=\", \"39+50=\", \"34+53=\", \"22+13=\"),\n    @(\"40+51=\", \"48-12=\", \"13+40=\", \"30+11=\", \"22+7=\"),\n    @(\"65+17=\", \"58-3=\", \"59+1=\", \"20+69=\", \"64-18=\"),\n    @(\"52+13=\", \"25+0=\", \"94-25=\", \"64-38=\", \"57-19=\"),\n    @(\"86-6=\", \"42-9=\", \"30+63=\", \"73-22=\", \"15-7=\"),\n    @(\"55-45=\", \"54-27=\", \"68+26=\", \"36+13=\", \"13-12=\"),\n    @(\"62+6=\", \"72+19=\", \"88-86=\", \"13+54=\", \"97-4=\"),\n    @(\"57-0=\", \"22+52=\", \"83+11=\", \"60-31=\", \"14-5=\"),\n    @(\"20+37=\", \"24+31=\", \"67-65=\", \"43-3=\", \"59-6=\"),\n)\n\nfor ($r = 0; $r -lt $newValues.Count; $r++) {\n    for ($c = 0; $c -lt $newValues[$r].Count; $c++) {\n        $t.Cell($r + 1, $c + 1).Range.Text = $newValues[$r][$c]\n    }\n}\n"}
